$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 46; everything from the old row 46 onward shifts down
# by one (old row 46 -> new row 47, ..., old row 60 -> new row 61).
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new data record.
$ws.Range("A46").Value = 11
$ws.Range("B46").Value = "Vega Monumental Concepción"
$ws.Range("C46").Value = "Bíobío"
$ws.Range("D46").Value = "2022-01-11"
$ws.Range("E46").Value = 8
$ws.Range("F46").Value = 100112001
$ws.Range("G46").Value = "Berenjena"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 310
$ws.Range("K46").Value = 6500
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = 6758
$ws.Range("N46").Value = "$/caja 60 unidades"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 113
$ws.Range("Q46").Value = 60
$ws.Range("R46").Value = "Hortaliza"
